$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (G1:I1), matching "only plotting k important features" ---
$ws.Range("G1").Value = "MSE_median"
$ws.Range("H1").Value = "MAE_median"
$ws.Range("I1").Value = "Dir_accuracy"

# Copy header formatting (bold, centered, bordered) from an existing header cell
# so the new header cells reuse the same style as the rest of row 1.
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)   # xlPasteFormats

# --- New metric values for DecisionTreeRegressor (row 2) ---
$ws.Range("G2").Value = 0.0004272416930457168
$ws.Range("H2").Value = 0.02066982395635263
$ws.Range("I2").Value = 0.3446808510638298

# --- New metric values for Naive baseline (row 3) ---
# Dir_accuracy (I3) is not computed for the Naive baseline, so it stays blank,
# mirroring the existing blank R^2 (F3) cell for that row.
$ws.Range("G3").Value = 0.001083194070471167
$ws.Range("H3").Value = 0.03291191380748265

$ws.Range("I3").IndentLevel = 0
$ws.Range("I3").ClearFormats()
